$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HTotalRevComparison")

# --- Header swap: B1 ("LATEST") <-> C1 ("PREVIOUS") ---
$ws.Range("B1").Value = "PREVIOUS"
$ws.Range("C1").Value = "LATEST"

# --- Column widths: column B widens from 9 to 10 (same as C:D) ---
$ws.Columns("B").ColumnWidth = 9.14

# --- Updated data values ---
$ws.Range("B2").Value = 2641.24
$ws.Range("C2").Value = 1152.02
$ws.Range("D2").Value = -1489.22

$ws.Range("B3").Value = 13278.34
$ws.Range("C3").Value = 13893.26
$ws.Range("D3").Value = 614.92

$ws.Range("B4").Value = 2024.62
$ws.Range("C4").Value = 1852.16
$ws.Range("D4").Value = -172.46

$ws.Range("B5").Value = 13231.05
$ws.Range("C5").Value = 14843.35
$ws.Range("D5").Value = 1612.3

$ws.Range("B6").Value = 4626.1
$ws.Range("C6").Value = 6237.46
$ws.Range("D6").Value = 1611.36

# --- D4 flips from "increase" (green) to "decrease" (red) styling ---
# --- D6 flips from "decrease" (red) to "increase" (green) styling ---
# Copy format (not value) from an existing cell with the right style so the
# shared style index is reused instead of minting a new one.
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
